$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Forandrad" (column C) date serial for all data rows 2-25: 46064 -> 46065
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 46065
}

# Reorder / update rows 7-25 (Beteckning, Datum, Markagare, Area) per source update

# Row 7
$ws.Cells.Item(7, 1).Value = "A 14674-2024"
$ws.Cells.Item(7, 2).Value = 45397
$ws.Cells.Item(7, 7).Value = 4.8

# Row 8
$ws.Cells.Item(8, 1).Value = "A 22256-2022"
$ws.Cells.Item(8, 2).Value = 44712
$ws.Cells.Item(8, 7).Value = 11.2

# Row 9
$ws.Cells.Item(9, 1).Value = "A 4574-2025"
$ws.Cells.Item(9, 2).Value = 45687
$ws.Cells.Item(9, 6).Value = "Kyrkan"
$ws.Cells.Item(9, 7).Value = 6.2

# Row 10
$ws.Cells.Item(10, 1).Value = "A 20239-2025"
$ws.Cells.Item(10, 2).Value = 45772
$ws.Cells.Item(10, 7).Value = 1.9

# Row 11
$ws.Cells.Item(11, 1).Value = "A 20054-2024"
$ws.Cells.Item(11, 2).Value = 45434
$ws.Cells.Item(11, 7).Value = 7.3

# Row 12
$ws.Cells.Item(12, 1).Value = "A 57394-2024"
$ws.Cells.Item(12, 2).Value = 45629.6907175926
$ws.Cells.Item(12, 7).Value = 0.5

# Row 13
$ws.Cells.Item(13, 2).Value = 44684

# Row 14
$ws.Cells.Item(14, 2).Value = 45629.68717592592

# Row 15
$ws.Cells.Item(15, 1).Value = "A 54557-2023"
$ws.Cells.Item(15, 2).Value = 45233.6346875
$ws.Cells.Item(15, 7).Value = 3.8

# Row 16
$ws.Cells.Item(16, 1).Value = "A 26708-2023"
$ws.Cells.Item(16, 2).Value = 45093
$ws.Cells.Item(16, 7).Value = 4.1

# Row 17
$ws.Cells.Item(17, 2).Value = 45295

# Row 18
$ws.Cells.Item(18, 1).Value = "A 59471-2024"
$ws.Cells.Item(18, 2).Value = 45638
$ws.Cells.Item(18, 7).Value = 1.8

# Row 19
$ws.Cells.Item(19, 1).Value = "A 46993-2025"
$ws.Cells.Item(19, 2).Value = 45929.54670138889
$ws.Cells.Item(19, 7).Value = 2.8

# Row 20
$ws.Cells.Item(20, 1).Value = "A 46998-2025"
$ws.Cells.Item(20, 2).Value = 45929.54851851852
$ws.Cells.Item(20, 7).Value = 0.9

# Row 21
$ws.Cells.Item(21, 1).Value = "A 53218-2023"
$ws.Cells.Item(21, 2).Value = 45229
$ws.Cells.Item(21, 7).Value = 5.4

# Row 22
$ws.Cells.Item(22, 1).Value = "A 63548-2025"
$ws.Cells.Item(22, 2).Value = 46013
$ws.Cells.Item(22, 7).Value = 0.9

# Row 23
$ws.Cells.Item(23, 1).Value = "A 6679-2026"
$ws.Cells.Item(23, 2).Value = 46056.60961805555
$ws.Cells.Item(23, 7).Value = 2.4

# Row 24
$ws.Cells.Item(24, 1).Value = "A 6684-2026"
$ws.Cells.Item(24, 2).Value = 46056.61989583333
$ws.Cells.Item(24, 7).Value = 8.199999999999999

# Row 25
$ws.Cells.Item(25, 1).Value = "A 59877-2025"
$ws.Cells.Item(25, 2).Value = 45993
$ws.Cells.Item(25, 6).ClearContents()
$ws.Cells.Item(25, 7).Value = 1
